# lead-bulk-template.xlsx: add "furnishingType" and "amenities" columns
# (changes in lead and broker)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at G:H ("flatType" stays at F, "areaKey"/"remark"
# that used to live in G/H shift right into I/J, carrying their old
# validation/formatting with them).
$ws.Columns("G:H").Insert()

# New header labels for the freshly inserted columns.
$ws.Range("G1").Value = "furnishingType"
$ws.Range("H1").Value = "amenities"

# Update the existing list validations that changed content.
$ws.Range("D2:D5000").Validation.Modify(3, 1, 1, '"10000-15000,15000-20000,20000-25000,25000-35000,35000-50000,50000-above"')
$ws.Range("E2:E5000").Validation.Modify(3, 1, 1, '"Standalone house,Apartment,Gated community,Independent house,Villa,PG / Co-living,Plot / Land,Anything is fine"')
$ws.Range("F2:F5000").Validation.Modify(3, 1, 1, '"1RK,1BHK,2BHK,3BHK,4BHK,Villa,Penthouse,Anything is fine"')

# Add the two new list validations for furnishingType / amenities.
$vG = $ws.Range("G2:G5000").Validation
$vG.Add(3, 1, 1, '"Fully Furnished,Semi Furnished,Unfurnished"')
$vG.ShowInput = $false
$vG.ShowError = $false

$vH = $ws.Range("H2:H5000").Validation
$vH.Add(3, 1, 1, '"Parking,Security,Power backup,Lift,Balcony"')
$vH.ShowInput = $false
$vH.ShowError = $false

# The areaKey validation (now on column I) keeps its original formula; just
# re-add it so it sorts after the new G/H validations like in the template.
$areaKeyFormula = '"Whitefield,Indiranagar,Koramangala,Bengaluru,Jayanagar,Banashankari,Basaveshwaranagar,Bheemanahalli,Bommanahalli,Chikkalasandra,Dasarahalli,Domlur,Electronic City,Frazer Town,Girinagar,Gokula,Gopalapuram,Hanumanthanagar,HBR Layout,Hebbal,Hoysala,HSR Layout,Ittamadu,JP Nagar,Jyothinagar,Kammanahalli,Kaval Byrasandra,Kodichikkanahalli,Kommadi,Kundalahalli,Lingrajapuram,Mahadevapura,Malleswaram,Marathahalli,Mathikere,Mico Layout,Mookambika,Nagavara,Nagawara,Nagarathpet,Nandini Layout,Nayandahalli,Old Airport Road,Peenya,Prithviraj Road,RMV Extension,Sadashivnagar,Sahakarnagar,Sanjaynagar,Sarjapur Road,Seshadripuram,Shantinagar,Shivaji Nagar,Soladevanahalli,Subramanyanagar"'
$ws.Range("I2:I5000").Validation.Delete()
$vI = $ws.Range("I2:I5000").Validation
$vI.Add(3, 1, 1, $areaKeyFormula)
$vI.ShowInput = $false
$vI.ShowError = $false

# Column widths (stored XML width = ColumnWidth + 5/6 in this engine).
$ws.Columns("C").ColumnWidth = 25 - 5/6
$ws.Columns("D").ColumnWidth = 18 - 5/6
$ws.Columns("E").ColumnWidth = 22 - 5/6
$ws.Columns("F").ColumnWidth = 18 - 5/6
$ws.Columns("G").ColumnWidth = 20 - 5/6
$ws.Columns("H").ColumnWidth = 30 - 5/6
$ws.Columns("I").ColumnWidth = 25 - 5/6
$ws.Columns("J").ColumnWidth = 30 - 5/6

Write-Host "Applied lead-bulk-template furnishingType/amenities column changes"
